# Adds data set additions
# Applies the "SCOPE_MODEL_DATA" sheet updates:
#  - row 4 / row 6 height adjustments
#  - J6 value swap ("randomly generated within range" -> "based on open_data returned type ")
#  - new "generation method" notes added to rows 16-18, plus a brand new row 19

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCOPE_MODEL_DATA")
$ws.Activate()

# --- row height tweaks -----------------------------------------------------
$ws.Rows.Item(4).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 34

# --- updated "generation method" note for the open-data altitude row -------
$ws.Range("J6").Value = "based on open_data returned type "

# --- new "generation method" notes for the image-based factor rows ---------
$ws.Range("J16").Value = "open_data"
$ws.Range("J17").Value = "from GTSBR {labelled}"
$ws.Range("J18").Value = "right"

# --- brand-new row 19 with its own generation-method note ------------------
$ws.Range("J19").Value = "10% left"
$ws.Rows.Item(19).RowHeight = 17

# --- refresh the on-screen view to match the authored selection/scroll -----
$ws.Range("C15").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
